$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 116, shifting the existing last data row (116) down to 117
$ws.Rows.Item(116).Insert()

# Populate the newly inserted row 116 with new data
$ws.Cells.Item(116, 1).Value = 3
$ws.Cells.Item(116, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(116, 3).Value = "Coquimbo"
$ws.Cells.Item(116, 4).Value = 44448
$ws.Cells.Item(116, 4).NumberFormat = $ws.Cells.Item(117, 4).NumberFormat
$ws.Cells.Item(116, 5).Value = 5
$ws.Cells.Item(116, 6).Value = 100112010
$ws.Cells.Item(116, 7).Value = "Achicoria"
$ws.Cells.Item(116, 8).Value = "Sin especificar"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 130
$ws.Cells.Item(116, 11).Value = 6000
$ws.Cells.Item(116, 12).Value = 6500
$ws.Cells.Item(116, 13).Value = 6269
$ws.Cells.Item(116, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(116, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(116, 16).Value = 392
$ws.Cells.Item(116, 17).Value = 16
$ws.Cells.Item(116, 18).Value = "Hortaliza"
